$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "concatenate"
$ws.Range("B3").Value = "Hello World"
